$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B17").Value = "Closed"
$ws.Range("B18").Value = "Scraped"
$ws.Range("B19").Value = "Scraped"
$ws.Range("B20").Value = "Scraped"
$ws.Range("B21").Value = "Scraped"
$ws.Range("B22").Value = "Failed"
